# Insert two new data rows (122-123) for a new weekly observation (D=44484),
# pushing all existing rows from 122 downward down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("122:123").Insert()

# Row 122 - "Primera" quality observation for the new date
$ws.Range("A122").Value = 8
$ws.Range("B122").Value = "Terminal La Palmera de La Serena"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 44484
$ws.Range("E122").Value = 4
$ws.Range("F122").Value = 100112009
$ws.Range("G122").Value = "Acelga"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 500
$ws.Range("L122").Value = 600
$ws.Range("M122").Value = 550
$ws.Range("N122").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O122").Value = "Provincia del Elquí"
$ws.Range("P122").Value = 275
$ws.Range("Q122").Value = 2
$ws.Range("R122").Value = "Hortaliza"

# Row 123 - "Segunda" quality observation for the same new date
$ws.Range("A123").Value = 8
$ws.Range("B123").Value = "Terminal La Palmera de La Serena"
$ws.Range("C123").Value = "Coquimbo"
$ws.Range("D123").Value = 44484
$ws.Range("E123").Value = 4
$ws.Range("F123").Value = 100112009
$ws.Range("G123").Value = "Acelga"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Segunda"
$ws.Range("J123").Value = 1480
$ws.Range("K123").Value = 400
$ws.Range("L123").Value = 450
$ws.Range("M123").Value = 425
$ws.Range("N123").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O123").Value = "Provincia del Elquí"
$ws.Range("P123").Value = 212
$ws.Range("Q123").Value = 2
$ws.Range("R123").Value = "Hortaliza"
